$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while guaranteeing it stays a TEXT string, even when the
# new value looks like a valid number (e.g. "322.36"). Excel would otherwise silently
# convert such assignments to a numeric cell. We flip the cell to text format ("@"),
# assign the value, then reset the style to Normal so no stray number-format/style
# index is left behind on the cell (done per-cell, immediately, to avoid multi-area
# range quirks).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Cell value updates (row order follows the sheet; rows 29 and 30 swap content) ---

$ws.Range("D2").Value = '42.408.96'
$ws.Range("E2").Value = '  +1.36%  '

$ws.Range("D3").Value = '2.291.86'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("E4").Value = '  +0.22%  '

Set-TextValue "D5" '322.36'
$ws.Range("E5").Value = '  +1.87%  '

Set-TextValue "D6" '104.22'
$ws.Range("E6").Value = '  +1.53%  '

Set-TextValue "D7" '0.629'
$ws.Range("E7").Value = '  +0.51%  '

$ws.Range("E8").Value = '  +0.16%  '

Set-TextValue "D9" '0.607'
$ws.Range("E9").Value = '  +0.81%  '

Set-TextValue "D10" '40.23'
$ws.Range("E10").Value = '  +3.55%  '

Set-TextValue "D11" '0.0906'
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("E12").Value = '  +3.07%  '

$ws.Range("E13").Value = '  +0.50%  '

Set-TextValue "D14" '0.966'
$ws.Range("E14").Value = '  +0.41%  '

Set-TextValue "D15" '15.23'
$ws.Range("E15").Value = '  -0.07%  '

$ws.Range("D16").Value = '2.639.97'
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("D17").Value = '2.284.57'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '42.396.74'
$ws.Range("E18").Value = '  +1.38%  '

Set-TextValue "D19" '7.44'
$ws.Range("E19").Value = '  -1.46%  '

$ws.Range("E20").Value = '  +0.15%  '

Set-TextValue "D21" '13.27'
$ws.Range("E21").Value = '  +34.04%  '

Set-TextValue "D22" '73.35'
$ws.Range("E22").Value = '  -0.49%  '

$ws.Range("E23").Value = '  +0.32%  '

Set-TextValue "D24" '269.54'
$ws.Range("E24").Value = '  -5.30%  '

$ws.Range("E25").Value = '  -1.63%  '

$ws.Range("E26").Value = '  -0.29%  '

Set-TextValue "D27" '10.88'
$ws.Range("E27").Value = '  +1.15%  '

$ws.Range("E28").Value = '  +2.18%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D29" '22.56'
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D30" '38.32'
$ws.Range("E30").Value = '  +10.48%  '

Set-TextValue "D31" '165.60'
$ws.Range("E31").Value = '  +1.67%  '

Set-TextValue "D32" '6.12'
$ws.Range("E32").Value = '  +4.76%  '

Set-TextValue "D33" '0.0882'
$ws.Range("E33").Value = '  +0.79%  '

Set-TextValue "D34" '0.132'
$ws.Range("E34").Value = '  +0.90%  '

$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("E36").Value = '  -13.17%  '

$ws.Range("E37").Value = '  +1.01%  '

$ws.Range("E38").Value = '  +1.72%  '

Set-TextValue "D39" '3.76'
$ws.Range("E39").Value = '  +4.68%  '

$ws.Range("E40").Value = '  -5.94%  '

Set-TextValue "D41" '1.54'
$ws.Range("E41").Value = '  +5.36%  '

Set-TextValue "D42" '69.90'
$ws.Range("E42").Value = '  +0.20%  '

Set-TextValue "D43" '95.81'
$ws.Range("E43").Value = '  -6.86%  '

$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("E46").Value = '  +3.60%  '

Set-TextValue "D47" '81.03'
$ws.Range("E47").Value = '  +5.47%  '

Set-TextValue "D48" '113.04'
$ws.Range("E48").Value = '  -2.02%  '

$ws.Range("E49").Value = '  -1.20%  '

Set-TextValue "D50" '5.26'
$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("D51").Value = '1.582.43'
$ws.Range("E51").Value = '  +3.04%  '
